$p = $ppt.ActivePresentation

# --- Slide 8 ("Coding Part"): "Lecture Folder: Lecture 1" -> "Lecture Folder: Lecture 2/Sequential/readPCAP.py"
$slide8 = $p.Slides.Item(8)
$shape8 = $slide8.Shapes.Item(2)
$tr8 = $shape8.TextFrame.TextRange
$old8 = "Lecture Folder: Lecture 1"
$new8 = "Lecture Folder: Lecture 2/Sequential/readPCAP.py"
$idx8 = $tr8.Text.IndexOf($old8)
if ($idx8 -ge 0) {
    $run8 = $tr8.Characters($idx8 + 1, $old8.Length)
    $run8.Text = $new8
}

# --- Slide 9 ("Extraction Field"): "Lecture Folder: Lecture 2" -> "Lecture Folder: Lecture 2/Sequential/readpcap_wireshark.py"
$slide9 = $p.Slides.Item(9)
$shape9 = $slide9.Shapes.Item(2)
$tr9 = $shape9.TextFrame.TextRange
$old9 = "Lecture Folder: Lecture 2"
$new9 = "Lecture Folder: Lecture 2/Sequential/readpcap_wireshark.py"
$idx9 = $tr9.Text.IndexOf($old9)
if ($idx9 -ge 0) {
    $run9 = $tr9.Characters($idx9 + 1, $old9.Length)
    $run9.Text = $new9
}
